$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "为AI加速而生：英特尔的至强，现在能跑200亿参数大模型了"
$ws.Range("B10").Value = "Weakly Supervised Open-Vocabulary Object Detection"
$ws.Range("B11").Value = "On Inference Stability for Diffusion Models"

$ws.Range("C9").Value = "https://www.jiqizhixin.com//articles/2023-12-20-3"
$ws.Range("C10").Value = "http://arxiv.org/abs/2312.12437v1"
$ws.Range("C11").Value = "http://arxiv.org/abs/2312.12431v1"

$ws.Range("D3").Value = "Wayve, a London-based autonomous driving technology company, is leading the way in AV 2.0, the next generation of self-driving cars. AV 2.0 focuses on comprehensive in-vehicle intelligence and uses generative AI to create and simulate novel driving scenarios. Wayve has developed GAIA-1, a generative world model, and LINGO-1, an AI model that enhances the learning and explainability of AI driving models. The company aims to improve the safety of autonomous vehicles and believes that embodied AI will play a definitive role in the future of the AI landscape."

$ws.Range("D6").Value = "OpenAI is conducting research to develop artificial general intelligence (AGI) that can solve human-level problems. They are focused on aligning powerful AI systems with human intentions and are using deep learning to build generative models for text, images, and audio. Their research has led to advancements in language processing, image generation, and audio processing. OpenAI has also developed tools like CLIP and DALL-E for connecting text and images. They are actively seeking talented individuals to join their team."

$ws.Range("D9").Value = "Intel has released its fifth-generation Xeon Scalable processors, which are designed to accelerate AI workloads. The new processors offer improved performance and efficiency, with up to 64 cores and AI acceleration capabilities in each core. They also feature Intel AMX, a dedicated matrix computing unit for AI acceleration. The processors have been validated by companies such as JD Cloud and Alibaba Cloud, which have reported significant performance improvements in AI workloads. Intel is also focusing on software support, providing tools like the oneAPI Deep Neural Network Library and OpenVINO toolkit to optimize and deploy AI models on its processors. The company is already working on the next generation of Xeon processors, which will have different core architectures to cater to different workloads."

$ws.Range("D10").Value = "The paper proposes a weakly supervised open-vocabulary object detection framework called WSOVOD. This framework extends traditional weakly supervised object detection (WSOD) to detect novel concepts and utilize diverse datasets with only image-level annotations. The authors explore three strategies: dataset-level feature adaptation, image-level salient object localization, and region-level vision-language alignment. Experimental results show that WSOVOD achieves state-of-the-art performance in both close-set object localization and detection tasks, and it also enables cross-dataset and open-vocabulary learning."

$ws.Range("D11").Value = "The paper discusses the issue of neglecting the correlation between timesteps in training methods for Denoising Probabilistic Models (DPMs), which limits their performance in generating images effectively. The authors propose a novel loss function that aims to reduce the estimation gap and enhance the sampling quality. Experimental results show that their method improves image generalization quality compared to several DPM baselines."

$ws.Range("D13").Value = "The RSNA 2023 conference has attracted over 40,000 registrants, according to Google News."
